$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the setpoint value in D7 (was 60, now 70).
# All downstream formula cells (rows 13, 14, 18, 19, 20) depend on this
# value and will recalculate automatically.
$ws.Range("D7").Value = 70

# Update the active cell selection to match the saved view state.
$null = $ws.Range("E8").Select()
